$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.252.74'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.268.27'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.14'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.38'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.54%  '
$ws.Range('E7').Value = '  -0.86%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.35'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.73%  '
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.87'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('D14').Value = '2.620.44'
$ws.Range('E14').Value = '  -0.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.77'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('D16').Value = '2.254.50'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '42.118.88'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.45'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.98%  '
$ws.Range('D20').Value = '0.0₃0908'
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.31'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.43'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.58'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.94'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.67'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.73'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.30%  '
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '162.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('E32').Value = '  -2.08%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  +3.18%  '
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.24'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.02'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.31'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.07'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.94%  '
$ws.Range('D44').Value = '1.943.03'
$ws.Range('E44').Value = '  -3.93%  '
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.00'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.34%  '
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.53'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '71.91'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '92.08'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  -1.34%  '
